$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting existing rows 28-33 down to 29-34.
$ws.Rows.Item(28).Insert()

# Copy the date-number style from D29 (old D28, shifted down) into new D28.
$ws.Range("D28").Value = 44754

$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112001
$ws.Range("G28").Value = "Berenjena"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 11000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 11500
$ws.Range("N28").Value = "$/caja 60 unidades"
$ws.Range("O28").Value = "Región de Arica y Parinacota"
$ws.Range("P28").Value = 192
$ws.Range("Q28").Value = 60
$ws.Range("R28").Value = "Hortaliza"

# Make sure D28 uses same style/number format as the other date cells in column D.
$ws.Range("D28").NumberFormat = $ws.Range("D29").NumberFormat
